$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "'-0.62%"
$ws.Range("D3").Value = "'27.02"
$ws.Range("E3").Value = "'0.51%"
$ws.Range("D4").Value = "'4.321"
$ws.Range("E4").Value = "'-9.53%"
$ws.Range("D5").Value = "'0.05888"
$ws.Range("E5").Value = "'-1.36%"
$ws.Range("D6").Value = "'6.633"
$ws.Range("E6").Value = "'-0.76%"
$ws.Range("D7").Value = "'0.8504"
$ws.Range("E7").Value = "'-3.06%"
$ws.Range("D8").Value = "'0.9373"
$ws.Range("E8").Value = "'-1.76%"
$ws.Range("D9").Value = "'0.1383"
$ws.Range("E9").Value = "'-2.33%"
$ws.Range("D10").Value = "'0.04668"
$ws.Range("E10").Value = "'29.28%"
$ws.Range("D11").Value = "'0.07074"
$ws.Range("E11").Value = "'-1.48%"
$ws.Range("D12").Value = "'0.03076"
$ws.Range("E12").Value = "'-2.14%"
$ws.Range("D13").Value = "'0.09117"
$ws.Range("E13").Value = "'-1.31%"
$ws.Range("D14").Value = "'0.001526"
$ws.Range("E14").Value = "'-0.96%"
$ws.Range("B15").Value = "One"
$ws.Range("C15").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D15").Value = "'0.0006052"
$ws.Range("E15").Value = "'-0.43%"
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").Value = "'0.006141"
$ws.Range("E16").Value = "'2.93%"
$ws.Range("B17").Value = "UpBots"
$ws.Range("C17").Value = "https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt"
$ws.Range("D17").Value = "'0.007492"
$ws.Range("E17").Value = "'4,919.57%"
$ws.Range("B18").Value = "LEO"
$ws.Range("C18").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D18").Value = "'3.491"
$ws.Range("E18").Value = "'0.14%"
$ws.Range("B19").Value = "GateToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D19").Value = "'3.169"
$ws.Range("E19").Value = "'-1.83%"
$ws.Range("B20").Value = "BTSEToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D20").Value = "'2.226"
$ws.Range("E20").Value = "'0.33%"
$ws.Range("B21").Value = "BitpandaEcosystemToken"
$ws.Range("C21").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D21").Value = "'0.3050"
$ws.Range("E21").Value = "'-2.70%"
$ws.Range("B22").Value = "ProBitToken"
$ws.Range("C22").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D22").Value = "'0.1270"
$ws.Range("E22").Value = "'-1.56%"
$ws.Range("B23").Value = "MCDex"
$ws.Range("C23").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D23").Value = "'3.921"
$ws.Range("E23").Value = "'10.95%"
$ws.Range("B24").Value = "CoinExToken"
$ws.Range("C24").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D24").Value = "'0.04252"
$ws.Range("E24").Value = "'0.85%"
$ws.Range("B25").Value = "BitKan"
$ws.Range("C25").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("D25").Value = "'0.001218"
$ws.Range("E25").Value = "'-0.29%"
$ws.Range("B26").Value = "HotbitToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("D26").Value = "'0.004282"
$ws.Range("E26").Value = "'-5.16%"
$ws.Range("B27").Value = "NitroEx"
$ws.Range("C27").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
$ws.Range("D27").Value = "'0.0001200"
$ws.Range("E27").Value = "'0.04%"
$ws.Range("D40").Value = "'0.03817"
$ws.Range("E40").Value = "'-0.82%"
$ws.Range("D41").Value = "'0.006211"
$ws.Range("E41").Value = "'3.71%"
$ws.Range("D42").Value = "'0.1101"
$ws.Range("E42").Value = "'-0.16%"
$ws.Range("D43").Value = "'0.002390"
$ws.Range("E43").Value = "'8.68%"
$ws.Range("D44").Value = "'0.01404"
$ws.Range("E44").Value = "'26.70%"
$ws.Range("D45").Value = "'0.00005367"
$ws.Range("E45").Value = "'-2.27%"
$ws.Range("E46").Value = "'0.07%"
$ws.Range("D47").Value = "'0.06590"
$ws.Range("E47").Value = "'-22.90%"
$ws.Range("D48").Value = "'0.2523"
$ws.Range("E48").Value = "'11,787.12%"
$ws.Range("D49").Value = "'0.00002101"
$ws.Range("E49").Value = "'0.07%"
$ws.Range("D50").Value = "'0.0002001"
$ws.Range("E50").Value = "'0.07%"
